# Regression Suite for CDS
# Applies the edits captured in the diff for
# TC10_CDS_phs002529_Platform_Illumina.xlsx:
#   1. The "Samples" query text (currently referenced by B3) loses its
#      Tumor / Analyte Type columns.
#   2. The shared-string table ends up with the Files query stored before
#      the (new, shorter) Samples query, which happens naturally because we
#      rewrite the Samples query text (forcing a new shared-string entry)
#      while the Files query text is untouched.
#   3. The view/selection on the sheet moves from C8 (topLeftCell B4) to C3
#      (topLeftCell B3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSamplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs002529' AND gi.platform = 'Illumina'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# Update the SamplesTab query (row 3, column B) to the new text that no
# longer projects Tumor / Analyte Type.
$ws.Range("B3").Value = $newSamplesQuery

# Move the view/selection to C3 (was C8), and scroll so row/col B3 is the
# top-left visible cell (was B4).
$win = $excel.ActiveWindow
$ws.Range("C3").Select()
$win.ScrollRow = 3
$win.ScrollColumn = 2
